$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cross Section Properties")

# Update stringer coordinates (D11/E11 lose their highlight style, D12/E12 keep it)
$ws.Range("D11").Formula = "=-1.6046754518"
$ws.Range("E11").Formula = "=-0.092597515"
$ws.Range("D12").Formula = "=-1.7254548192"
$ws.Range("E12").Formula = "=-0.1248053464"

# D11/E11 previously had style index 5 (applyFill, no fill) - clear that override so it
# reverts to the default "Normal" style, matching the target state.
$ws.Range("D11:E11").Style = "Normal"

# Update the active selection to match the authored state
$ws.Range("D13").Select()

$wb.Application.Calculate()
